# Regional- Appendix.xlsx — "Had to update data with complete dataset back to
# 1980 then reprocess."
#
# Eleven data points on the "Regional for Mapping" sheet are corrected with
# reprocessed values. All of the source cells are stored as text (not
# numbers), so a leading apostrophe is used to force the replacement values
# to stay text as well (mirrors typing '0.5 etc. directly into Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regional for Mapping")

$ws.Range("AB2").Value  = "'0.5"
$ws.Range("X3").Value   = "'2.2842"
$ws.Range("AD6").Value  = "'41"
$ws.Range("AD9").Value  = "'18"
$ws.Range("AD10").Value = "'9"
$ws.Range("AD11").Value = "'19"
$ws.Range("AD12").Value = "'10"
$ws.Range("X13").Value  = "'3.7858"
$ws.Range("X21").Value  = "'4.8775"
$ws.Range("X22").Value  = "'3.0"
$ws.Range("X23").Value  = "'1.1246"

# The "Regional Exhibit" sheet's first and sixth data columns were also
# widened slightly (likely to fit the refreshed figures/labels).
$ws2 = $wb.Worksheets.Item("Regional Exhibit")
$ws2.Columns.Item(1).ColumnWidth = 11.25
$ws2.Columns.Item(6).ColumnWidth = 15.25
